# issue #5: stock data from json to db
# Adds a "category" column (value "normal") right after "property_category",
# and two trailing columns "source_file" (value "tmp9bfb1") and "index"
# (the original row's index number) to the 股票 (stock) worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column before the existing "date" column (column I) so the
# old I/J/K (date/legislator_name/legislator_id) shift one column right and
# make room for the new "category" column.
$ws.Columns.Item(9).Insert()

# --- Header row (row 1) ---
$ws.Cells.Item(1, 9).Value2 = "category"
$ws.Cells.Item(1, 13).Value2 = "source_file"
$ws.Cells.Item(1, 14).Value2 = "index"

# Match the bold/centered/bordered header style used by the other header
# cells (B1:L1) for the two newly appended header cells (M1, N1).
foreach ($col in 13, 14) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# --- Data rows (row 2..6) ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $origIndex = $ws.Cells.Item($r, 1).Value2

    $ws.Cells.Item($r, 9).Value2 = "normal"
    $ws.Cells.Item($r, 13).Value2 = "tmp9bfb1"
    $ws.Cells.Item($r, 14).Value2 = $origIndex
}
